$wb = $excel.ActiveWorkbook

# --- Update the "Date" value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# --- Update the Display column (C) capitalization on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C5").Value = "Prefer not to Answer"
$concepts.Range("C6").Value = "Not Applicable"
$concepts.Range("C7").Value = "Missing - Unknown"
$concepts.Range("C8").Value = "Missing - Not Collected"
$concepts.Range("C9").Value = "Missing - Not Provided"
$concepts.Range("C10").Value = "Missing - Restricted Access"
